$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Br"
$ws.Range("A4").Value = "CaCO3"
$ws.Range("A4").Select()
